# Split "Ministry Course Code and Level" (e.g. "ENST 12") into two
# separate columns: "Ministry Course Code" (e.g. "ENST") and
# "Ministry Course Level" (e.g. 12) - Summer Reporting File Spec change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the existing "Ministry Course Code and Level"
# column (G), pushing that column (and its data, e.g. "ENST 12") one
# column to the right, to H.
$ws.Columns("G").Insert() | Out-Null

# The freshly-inserted column G picks up the formatting of its left
# neighbour (F); match it up with the rest of the data columns (same
# left-aligned / general-format style used by the other data cells)
# by copying the format from a representative data cell.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("G2:G4").PasteSpecial(-4122) | Out-Null

# Re-label the (now split) course-code column and populate it with just
# the course code.
$ws.Range("G1").Value = "Ministry Course Code"
$ws.Range("G2:G4").Value = "ENST"

# The old "Ministry Course Code and Level" column (now shifted to H,
# still holding "ENST 12") becomes the course-level column.
$ws.Range("H1").Value = "Ministry Course Level"
$ws.Range("H2:H4").Value = 12

# Update the selection to match the new "split" columns.
$ws.Range("G1:H1048576").Select() | Out-Null
